# Task_Utilisation.xlsx refresh — pulls in the latest FreeRTOS task-utilisation
# snapshot (stack high-water marks + per-task runtime percentages) captured
# from the device, and moves the on-screen selection to AA9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Utilisation")

# --- Task 1: defaultTask (row 2) ---------------------------------------
# Runtime %
$ws.Range("H2").Value = "0x3122 (1.9%)"

# --- Task 2: SPIDeviceManager (row 7) -----------------------------------
# Runtime %
$ws.Range("H7").Value = "0xe86f (9.2%)"

# --- Task 3: USARTDeviceManager (rows 12-16) ----------------------------
# Stack Usage snapshot (Q/R/S/T columns are computed from this)
$ws.Range("F12").Value = "452 B / 1.48 kB"
# Runtime %
$ws.Range("H12").Value = "0x4ed2 (3.1%)"
# Stack High Water Mark
$ws.Range("C16").Value = "0x20002734"

# --- Task 4: I2CDeviceManager (row 17) ----------------------------------
# Runtime %
$ws.Range("H17").Value = "0x7963 (4.8%)"

# --- Task 5: ADCDeviceManager (row 22) ----------------------------------
# Runtime %
$ws.Range("H22").Value = "0x680a (4.1%)"

# --- Task 6: FANMotorManager (row 27) -----------------------------------
# Runtime %
$ws.Range("H27").Value = "0xf (0.0%)"

# --- Task 7: DACDeviceManager (row 32) ----------------------------------
# Runtime %
$ws.Range("H32").Value = "0x96a (0.4%)"

# --- Task 8: STPMotorManager (row 37) -----------------------------------
# Runtime %
$ws.Range("H37").Value = "0x112 (0.0%)"

# --- Task 9: IDLE (rows 42-46) ------------------------------------------
# Stack Usage snapshot
$ws.Range("F42").Value = "88 B / 248 B"
# Runtime %
$ws.Range("H42").Value = "0x78702 (76.4%)"
# Stack High Water Mark
$ws.Range("C46").Value = "0x20003ca0"

# Move/refresh the on-screen selection to match the author's last cursor
# position when they saved.
$ws.Range("AA9").Select()
